# Actualización automática 2025-07-31 08:55:10
$wb = $excel.ActiveWorkbook

$wsVentasPorGrupo   = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentaMensual     = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumplimiento     = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# ----- Sheet "VENTAS POR GRUPO" -----
$wsVentasPorGrupo.Range("C4").Value = 2073.6
$wsVentasPorGrupo.Range("M4").Value = 9199.01

$wsVentasPorGrupo.Range("M5").Value = 4930.08

$wsVentasPorGrupo.Range("D16").Value = 457.92
$wsVentasPorGrupo.Range("L16").Value = 855.36
$wsVentasPorGrupo.Range("M16").Value = 5590.26
$wsVentasPorGrupo.Range("N16").Value = 262.85

$wsVentasPorGrupo.Range("C25").Value = 2052.86
$wsVentasPorGrupo.Range("L25").Value = 1088.66
$wsVentasPorGrupo.Range("M25").Value = 9445.120000000001

$wsVentasPorGrupo.Range("I54").Value = 26.1

$wsVentasPorGrupo.Range("D56").Value = "1 de 54"
$wsVentasPorGrupo.Range("I56").Value = "6 de 54"
$wsVentasPorGrupo.Range("L56").Value = "9 de 54"
$wsVentasPorGrupo.Range("N56").Value = "1 de 54"

# ----- Sheet "VENTA MENSUAL" -----
$wsVentaMensual.Range("F4").Value = 11837.85
$wsVentaMensual.Range("F5").Value = 7761.47
$wsVentaMensual.Range("F16").Value = 7372.95
$wsVentaMensual.Range("F25").Value = 13349.81
$wsVentaMensual.Range("F54").Value = 820.98
$wsVentaMensual.Range("F56").Value = 81130.25

# ----- Sheet "CUMPLIMIENTO MENSUAL" -----
$wsCumplimiento.Range("D2").Value = 5681.66
$wsCumplimiento.Range("E2").Value = 4288.68304517915
$wsCumplimiento.Range("F2").Value = 0.5698560194222395

$wsCumplimiento.Range("D3").Value = 457.92
$wsCumplimiento.Range("E3").Value = 26999.0876
$wsCumplimiento.Range("F3").Value = 0.01667770962776002

$wsCumplimiento.Range("D8").Value = 692.86
$wsCumplimiento.Range("E8").Value = 307.14
$wsCumplimiento.Range("F8").Value = 0.69286

$wsCumplimiento.Range("D15").Value = 7978.29
$wsCumplimiento.Range("E15").Value = 5521.71
$wsCumplimiento.Range("F15").Value = 0.5909844444444444

$wsCumplimiento.Range("D16").Value = 62700.73
$wsCumplimiento.Range("E16").Value = -10874.27
$wsCumplimiento.Range("F16").Value = 1.209820813538104

$wsCumplimiento.Range("D17").Value = 262.85
$wsCumplimiento.Range("E17").Value = 421.15
$wsCumplimiento.Range("F17").Value = 0.3842836257309942

$wsCumplimiento.Range("D19").Value = 81130.25000000001
$wsCumplimiento.Range("E19").Value = 32576.20064517915
$wsCumplimiento.Range("F19").Value = 0.7135061339058667
